$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 785, shifting existing rows 785:851 down to 786:852.
$ws.Rows.Item(785).Insert()

# Populate the newly inserted row 785 with the new weekly price record.
$row = 785
$ws.Cells.Item($row, 1).Value  = 10
$ws.Cells.Item($row, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value  = "La Araucanía"
$ws.Cells.Item($row, 4).Value  = 45223
$ws.Cells.Item($row, 5).Value  = 9
$ws.Cells.Item($row, 6).Value  = 100112027
$ws.Cells.Item($row, 7).Value  = "Melón"
$ws.Cells.Item($row, 8).Value  = "Tuna"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 2000
$ws.Cells.Item($row, 12).Value = 2000
$ws.Cells.Item($row, 13).Value = 2000
$ws.Cells.Item($row, 14).Value = "`$/unidad"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 2000
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
